# Plantilla Lista de Tareas de la 6ta Iteración - iteration template update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 12: task status changed from "Por iniciar" to "Hecho", and 2 hours
# consumed on day 6 (column W) are now logged.
$ws.Range("F12").Value = "Hecho"
$ws.Range("W12").Value = 2

# Row 14: the task entry ("Actualizar diagrama de paquetes.", Mario,
# Hecho, 2h) is removed entirely, clearing the task description,
# responsible person, status and estimated hours, as well as the 2 hours
# that had been logged against it on day 6 (column W).
$ws.Range("D14:G14").ClearContents()
$ws.Range("W14").ClearContents()

# Restore the blank-row look of E14 (status column) to match the style
# used by other empty rows instead of the colored "Hecho" status fill.
$ws.Range("E15").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to reflect where the editor was last working.
$ws.Range("F14").Select() | Out-Null
